# Braga-demo: remove names of concepts and relations, drop the Intersections sheet.

$wb = $excel.ActiveWorkbook

# 1. Remove the "Intersections" worksheet entirely.
$wsIntersections = $wb.Worksheets.Item("Intersections")
$wsIntersections.Delete() | Out-Null

# 2. Concepts sheet: drop the "name" and "isa" columns (B and C),
#    keeping only the identifier column (A). The leftover "name"/"Identifier"
#    labels are relocated to E8/E9.
$wsConcepts = $wb.Worksheets.Item("Concepts")
$wsConcepts.Columns("B:C").Delete() | Out-Null
$wsConcepts.Range("E8").Value = "name"
$wsConcepts.Range("E9").Value = "Identifier"

# 3. Relations sheet: drop the "name" column (B), keeping identifier (A),
#    source (now B) and target (now C).
$wsRelations = $wb.Worksheets.Item("Relations")
$wsRelations.Columns("B").Delete() | Out-Null

# 4. Restore/update the selections shown in each sheet view.
$wsConcepts.Activate()
$wsConcepts.Range("B1:C5").Select() | Out-Null

$wsRelations.Activate()
$wsRelations.Range("A6:XFD14").Select() | Out-Null

$wsCompositions = $wb.Worksheets.Item("Compositions")
$wsCompositions.Activate()
$wsCompositions.Range("C21").Select() | Out-Null
